$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (diem_tichluy / diem_renluyen swap columns H/I stays same text,
#     but the old "xet_hocbong" column J is being dropped) ---
$ws.Range("H1").Value = "diem_tichluy"
$ws.Range("I1").Value = "diem_renluyen"

# Remove the now-unused "xet_hocbong" column entirely (was column J)
$ws.Columns("J").Delete()

# Make sure the phone-number column keeps its leading zeros as text,
# matching the existing style used on G2/G3 in the original sheet.
$ws.Range("G2:G6").NumberFormat = "@"

# --- Row 2: Nguyen Thi Hai ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 111
$ws.Range("C2").Value = "Nguyen Thi"
$ws.Range("D2").Value = "Hai"
$ws.Range("E2").Value = "haint"
$ws.Range("F2").Value = "haint@vnuf.edu.vn"
$ws.Range("G2").Value = "0988111111"
$ws.Range("H2").Value = 2.2
$ws.Range("I2").Value = 10

# --- Row 3: Tran Van Thanh ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 222
$ws.Range("C3").Value = "Tran Van"
$ws.Range("D3").Value = "Thanh"
$ws.Range("E3").Value = "thanhtv"
$ws.Range("F3").Value = "thanhtv@vnuf.edu.vn"
$ws.Range("G3").Value = "0978222222"
$ws.Range("H3").Value = 2.8
$ws.Range("I3").Value = 7

# --- Row 4: Dang Thai Chau ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 333
$ws.Range("C4").Value = "Dang Thai"
$ws.Range("D4").Value = "Chau"
$ws.Range("E4").Value = "chaudt"
$ws.Range("F4").Value = "chaudt@vnuf.edu.vn"
$ws.Range("G4").Value = "0966333333"
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 6.5

# --- Row 5: Hoang Hai Nam ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 444
$ws.Range("C5").Value = "Hoang Hai"
$ws.Range("D5").Value = "Nam"
$ws.Range("E5").Value = "namhn"
$ws.Range("F5").Value = "namhn@vnuf.edu.vn"
$ws.Range("G5").Value = "0912444444"
$ws.Range("H5").Value = 1.6
$ws.Range("I5").Value = 8.8

# --- Row 6: Dinh Manh Thang ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 555
$ws.Range("C6").Value = "Dinh Manh"
$ws.Range("D6").Value = "Thang"
$ws.Range("E6").Value = "thangdm"
$ws.Range("F6").Value = "thangdm@vnuf.edu.vn"
$ws.Range("G6").Value = "0923555555"
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 5

# Move the active selection the way the saved file shows (cursor parked
# below the data, on column G).
$ws.Range("G7").Select()
